$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize font colour on a handful of existing cells (theme colour -> explicit black)
$ws.Range("C1").Font.Color = 0
$ws.Range("C2").Font.Color = 0
$ws.Range("B3").Font.Color = 0
$ws.Range("B4").Font.Color = 0
$ws.Range("C5").Font.Color = 0

# B5 becomes left aligned (previously "general")
$ws.Range("B5").HorizontalAlignment = -4131

# New object: "caixa de som" (Alexa)
$ws.Range("A6:E6").Style = "Normal"
$ws.Range("A6").Value = "Alexa"
$ws.Range("B6").Value = $false
$ws.Range("C6").Value = 50
$ws.Range("D6").Value = "Leo Santana"
$ws.Range("E6").Value = $true
$ws.Rows(6).AutoFit()

# Row 1 grew to match the standard row height used elsewhere in the sheet
$ws.Rows(1).RowHeight = 19.5
